# Moved to U.S. version 2.1.1 as baseline
# Update the CPI workbook: roll the "current" year forward from 2018 to 2019,
# refresh the BLS source link, and append the 2019 CPI-U data row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("Data")

# --- "About" sheet: bump the headline year -------------------------------
$ws1.Range("B4").Value = 2019

# --- "Data" sheet: append the new 2019 row (row 57) -----------------------
# NOTE: set the new "2019....." label text before the new hyperlink text so
# the shared-string table grows in the same order as the reference edit
# (index 72 = "2019....." label, index 73 = refreshed URL string).
$ws2.Range("A57").Value = "2019.............................................................................     ."
$ws2.Range("B57").Value = 254.41200000000001
$ws2.Range("C57").Value = 256.90300000000002
$ws2.Range("D57").Value = 255.65700000000001
$ws2.Range("E57").Value = 2.2999999999999998
$ws2.Range("F57").Value = 1.8

# Extend the shared "multiply by to get 2012 dollars" formula down to row 57.
$ws2.Range("G57").Formula = "=`$D`$50/D57"
$ws2.Range("G57").NumberFormat = "0.000"

# --- "About" sheet: refresh the BLS source URL text -----------------------
$ws1.Range("B6").Value = "https://www.bls.gov/cpi/tables/supplemental-files/historical-cpi-u-201912.pdf"
